# Apply updated cryptocurrency price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.204.09"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "3.712.15"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'236.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  +1.13%  "
$ws.Range("D7").Value = "'657.72"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.431"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  -2.41%  "
$ws.Range("D11").Value = "3.709.60"
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").Value = "'0.0000319"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +17.51%  "
$ws.Range("D13").Value = "'44.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.33%  "
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("D15").Value = "'6.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.81%  "
$ws.Range("D16").Value = "4.404.76"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("D17").Value = "96.759.63"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "'8.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").Value = "3.722.43"
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("E20").Value = "  +1.80%  "
$ws.Range("D21").Value = "'18.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.64%  "
$ws.Range("D22").Value = "'0.505"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.39%  "
$ws.Range("D23").Value = "'524.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("D24").Value = "'3.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.51%  "
$ws.Range("D25").Value = "'0.0000223"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.31%  "
$ws.Range("D26").Value = "'6.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.96%  "
$ws.Range("D27").Value = "'106.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.22%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "3.912.97"
$ws.Range("E28").Value = "  +0.90%  "
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").Value = "'0.190"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +13.00%  "
$ws.Range("D30").Value = "'13.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("D31").Value = "'12.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.54%  "
$ws.Range("D32").Value = "'3.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.71%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E34").Value = "  +3.11%  "
$ws.Range("D35").Value = "'1.83"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.96%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").Value = "'32.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("D38").Value = "'637.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.76%  "
$ws.Range("E39").Value = "  -1.54%  "
$ws.Range("D40").Value = "'8.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "'0.166"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.01%  "
$ws.Range("D43").Value = "'40.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.10%  "
$ws.Range("D44").Value = "'6.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.38%  "
$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D45").Value = "'2.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "'0.483"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +12.53%  "
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("E48").Value = "  -1.41%  "
$ws.Range("D49").Value = "'2.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.72%  "
$ws.Range("D50").Value = "'23.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("E51").Value = "  -0.90%  "
